# Insert a new weekly record row at row 229 (pushing existing rows 229:242
# down to 230:243), then populate the new row with the latest weekly data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(229).Insert()

$ws.Range("A229").Value = 9
$ws.Range("B229").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C229").Value = "Metropolitana"
$ws.Range("D229").Value = 44516
$ws.Range("E229").Value = 13
$ws.Range("F229").Value = 100112044
$ws.Range("G229").Value = "Perejil"
$ws.Range("H229").Value = "Sin especificar"
$ws.Range("I229").Value = "Primera"
$ws.Range("J229").Value = 106
$ws.Range("K229").Value = 10000
$ws.Range("L229").Value = 12000
$ws.Range("M229").Value = 11000
$ws.Range("N229").Value = "`$/docena de atados"
$ws.Range("O229").Value = "Región Metropolitana"
$ws.Range("P229").Value = 3667
$ws.Range("Q229").Value = 3
$ws.Range("R229").Value = "Hortaliza"
